$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price / 1h-volume data, plus the VeChain / ApeXProtocol row swap.
# Cells are stored as text in the source sheet (t="inlineStr"), so each cell's number
# format is briefly switched to text ("@") while the new value is written, then restored
# to "General" so Excel does not auto-convert numeric-looking strings (e.g. "1.00") to numbers.

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "66.206.70"
$cell.NumberFormat = "General"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -0.66%  "
$cell.NumberFormat = "General"

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.516.17"
$cell.NumberFormat = "General"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  +1.25%  "
$cell.NumberFormat = "General"

# Row 4
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.NumberFormat = "General"
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  +0.04%  "
$cell.NumberFormat = "General"

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "573.78"
$cell.NumberFormat = "General"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  +5.12%  "
$cell.NumberFormat = "General"

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "179.09"
$cell.NumberFormat = "General"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -4.49%  "
$cell.NumberFormat = "General"

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.637"
$cell.NumberFormat = "General"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  +5.61%  "
$cell.NumberFormat = "General"

# Row 8
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -0.11%  "
$cell.NumberFormat = "General"

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.636"
$cell.NumberFormat = "General"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  +1.52%  "
$cell.NumberFormat = "General"

# Row 10
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  +4.66%  "
$cell.NumberFormat = "General"

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "55.67"
$cell.NumberFormat = "General"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  +2.25%  "
$cell.NumberFormat = "General"

# Row 12
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +3.47%  "
$cell.NumberFormat = "General"

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "9.31"
$cell.NumberFormat = "General"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +0.15%  "
$cell.NumberFormat = "General"

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.074.22"
$cell.NumberFormat = "General"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  +1.33%  "
$cell.NumberFormat = "General"

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.512.53"
$cell.NumberFormat = "General"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  +1.18%  "
$cell.NumberFormat = "General"

# Row 16
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  +0.35%  "
$cell.NumberFormat = "General"

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "18.37"
$cell.NumberFormat = "General"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  +1.93%  "
$cell.NumberFormat = "General"

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "66.146.26"
$cell.NumberFormat = "General"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -0.73%  "
$cell.NumberFormat = "General"

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.04"
$cell.NumberFormat = "General"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  +2.96%  "
$cell.NumberFormat = "General"

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.NumberFormat = "General"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  +1.97%  "
$cell.NumberFormat = "General"

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "415.54"
$cell.NumberFormat = "General"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -0.33%  "
$cell.NumberFormat = "General"

# Row 22
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  +8.20%  "
$cell.NumberFormat = "General"

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.30"
$cell.NumberFormat = "General"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  +3.62%  "
$cell.NumberFormat = "General"

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "85.79"
$cell.NumberFormat = "General"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +2.01%  "
$cell.NumberFormat = "General"

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "13.05"
$cell.NumberFormat = "General"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  +10.58%  "
$cell.NumberFormat = "General"

# Row 26
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -0.94%  "
$cell.NumberFormat = "General"

# Row 27
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -0.73%  "
$cell.NumberFormat = "General"

# Row 28
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  +4.29%  "
$cell.NumberFormat = "General"

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "30.52"
$cell.NumberFormat = "General"

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "639.44"
$cell.NumberFormat = "General"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -2.21%  "
$cell.NumberFormat = "General"

# Row 31
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -1.40%  "
$cell.NumberFormat = "General"

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "11.71"
$cell.NumberFormat = "General"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +0.56%  "
$cell.NumberFormat = "General"

# Row 33
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  +1.31%  "
$cell.NumberFormat = "General"

# Row 34
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +14.18%  "
$cell.NumberFormat = "General"

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "59.61"
$cell.NumberFormat = "General"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +1.04%  "
$cell.NumberFormat = "General"

# Row 36
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.NumberFormat = "General"

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0799"
$cell.NumberFormat = "General"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -1.31%  "
$cell.NumberFormat = "General"

# Row 38
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "37.28"
$cell.NumberFormat = "General"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -2.26%  "
$cell.NumberFormat = "General"

# Row 39
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -1.57%  "
$cell.NumberFormat = "General"

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "3.258.73"
$cell.NumberFormat = "General"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +9.14%  "
$cell.NumberFormat = "General"

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "3.43"
$cell.NumberFormat = "General"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  +2.80%  "
$cell.NumberFormat = "General"

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.NumberFormat = "General"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +0.17%  "
$cell.NumberFormat = "General"

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.91"
$cell.NumberFormat = "General"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  +1.01%  "
$cell.NumberFormat = "General"

# Row 44
$cell = $ws.Range("B44")
$cell.NumberFormat = "@"
$cell.Value = "ApeXProtocol"
$cell.NumberFormat = "General"
$cell = $ws.Range("C44")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$cell.NumberFormat = "General"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "3.30"
$cell.NumberFormat = "General"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -3.86%  "
$cell.NumberFormat = "General"

# Row 45
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -3.23%  "
$cell.NumberFormat = "General"

# Row 46
$cell = $ws.Range("B46")
$cell.NumberFormat = "@"
$cell.Value = "VeChain"
$cell.NumberFormat = "General"
$cell = $ws.Range("C46")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell.NumberFormat = "General"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0418"
$cell.NumberFormat = "General"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  +0.96%  "
$cell.NumberFormat = "General"

# Row 47
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +1.24%  "
$cell.NumberFormat = "General"

# Row 48
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  +2.95%  "
$cell.NumberFormat = "General"

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "8.69"
$cell.NumberFormat = "General"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -0.79%  "
$cell.NumberFormat = "General"

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "139.68"
$cell.NumberFormat = "General"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  +0.26%  "
$cell.NumberFormat = "General"

# Row 51
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -0.60%  "
$cell.NumberFormat = "General"
